$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.972.24"
$ws.Range("E2").Value = "  +2.18%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.846.81"
$ws.Range("E3").Value = "  +2.32%  "

# Row 4
$ws.Range("E4").Value = "  -0.06%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "310.66"
$ws.Range("E5").Value = "  +1.28%  "

# Row 6
$ws.Range("E6").Value = "  -0.10%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4670"
$ws.Range("E7").Value = "  +3.35%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3633"
$ws.Range("E8").Value = "  +1.09%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07179"
$ws.Range("E9").Value = "  +1.56%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9263"
$ws.Range("E10").Value = "  +4.24%  "

# Row 11
$ws.Range("E11").Value = "  +1.12%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07701"
$ws.Range("E12").Value = "  -1.40%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.849.37"
$ws.Range("E13").Value = "  +2.39%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.303"
$ws.Range("E14").Value = "  +0.51%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.407"
$ws.Range("E15").Value = "  +1.70%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "88.46"
$ws.Range("E16").Value = "  +3.79%  "

# Row 17
$ws.Range("E17").Value = "  -0.09%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008605"
$ws.Range("E18").Value = "  +1.25%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "26.996.06"
$ws.Range("E20").Value = "  +2.06%  "

# Row 21
$ws.Range("E21").Value = "  +1.84%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.038"
$ws.Range("E22").Value = "  +1.62%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.64"
$ws.Range("E23").Value = "  +1.37%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.930"
$ws.Range("E24").Value = "  -1.54%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "152.56"
$ws.Range("E25").Value = "  +0.13%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "18.22"
$ws.Range("E26").Value = "  +2.41%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.056"
$ws.Range("E27").Value = "  -0.52%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "114.06"
$ws.Range("E28").Value = "  +1.75%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.932"
$ws.Range("E29").Value = "  +1.85%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.08863"
$ws.Range("E30").Value = "  +2.00%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.181"
$ws.Range("E31").Value = "  +3.31%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.854"
$ws.Range("E32").Value = "  +2.04%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.179"
$ws.Range("E33").Value = "  +6.64%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7478"
$ws.Range("E34").Value = "  +3.38%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.476"
$ws.Range("E35").Value = "  +0.51%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.088"
$ws.Range("E36").Value = "  +1.02%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.987"
$ws.Range("E37").Value = "  +2.58%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01941"
$ws.Range("E38").Value = "  +0.43%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05171"
$ws.Range("E39").Value = "  +1.27%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5165"
$ws.Range("E40").Value = "  +2.33%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.889"
$ws.Range("E41").Value = "  +1.65%  "

# Row 42
$ws.Range("E42").Value = "  +0.05%  "

# Row 43
$ws.Range("B43").Value = "Aptos"
$ws.Range("C43").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.186"
$ws.Range("E43").Value = "  +2.15%  "

# Row 44
$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "10.57"
$ws.Range("E44").Value = "  +5.65%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.4718"
$ws.Range("E45").Value = "  +1.22%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.008"
$ws.Range("E46").Value = "  -0.03%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "100.59"
$ws.Range("E47").Value = "  +0.60%  "

# Row 48
$ws.Range("E48").Value = "  +1.78%  "

# Row 49
$ws.Range("E49").Value = "  +1.24%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "64.40"
$ws.Range("E50").Value = "  +1.40%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "36.14"
$ws.Range("E51").Value = "  +0.16%  "
